$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark completed tasks (column C = "x") for rows 7 and 9
$ws.Range("C7").Value = "x"
$ws.Range("C9").Value = "x"

# Add a new bug entry in row 19
$ws.Range("B19").Value = "fix background notification bug caused by launcing application w/o internet and then logging in with internet"

# Update selection to the newly added cell
$ws.Range("B19").Select()
